{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same net changes as the target diff:\n//   1. \"Constructors and Classes will be in Upper CamelCase\"\n//        -> \"Constructors and Classes will however be in Upper CamelCase\"\n//   2. \"...in the JavaScript code as it can cause errors\"\n//        -> \"...in the JavaScript code as it can cause escape errors\"\n//   3. Append a new (empty) table row numbered \"14\" at the end of the table.\n//\n// (The diff also contains several hunks where Word merely re-split or\n// re-merged adjacent runs with no visible text change -- those are not\n// meaningful content edits and are not reproduced here.)\n\nconst body = context.document.body;\n\n// --- 1. Insert \"however\" into the Constructors/Classes sentence ---------\nconst ctorResults = body.search(\n  \"Constructors and Classes will be in Upper CamelCase\",\n  { matchCase: true }\n);\nctorResults.load(\"items\");\nawait context.sync();\n\nif (ctorResults.items.length > 0) {\n  ctorResults.items[0].insertText(\n    \"Constructors and Classes will however be in Upper CamelCase\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 2. Insert \"escape\" before \"errors\" in the quotation marks rule -----\nconst quoteResults = body.search(\n  \"Use only single quotation marks for example (\\u2018Hello World\\u2019) in the JavaScript code as it can cause errors\",\n  { matchCase: true }\n);\nquoteResults.load(\"items\");\nawait context.sync();\n\nif (quoteResults.items.length > 0) {\n  quoteResults.items[0].insertText(\n    \"Use only single quotation marks for example (\\u2018Hello World\\u2019) in the JavaScript code as it can cause escape errors\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 3. Append a new, mostly-empty row (convention #14) to the table ----\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.addRows(Word.InsertLocation.end, 1, [[\"14\", \"\", \"\"]]);\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the same net changes as the target diff:\n#   1. \"Constructors and Classes will be in Upper CamelCase\"\n#        -> \"Constructors and Classes will however be in Upper CamelCase\"\n#   2. \"...in the JavaScript code as it can cause errors\"\n#        -> \"...in the JavaScript code as it can cause escape errors\"\n#   3. Append a new (empty) table row numbered \"14\" at the end of the table.\n#\n# (The diff also contains several hunks where Word merely re-split or\n# re-merged adjacent runs with no visible text change -- those are not\n# meaningful content edits and are not reproduced here.)\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert \"however\" into the Constructors/Classes sentence ---------\n$find1 = $d.Content.Find\n$find1.Text = \"Constructors and Classes will be in Upper CamelCase\"\n$find1.Replacement.Text = \"Constructors and Classes will however be in Upper CamelCase\"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# --- 2. Insert \"escape\" before \"errors\" in the quotation marks rule -----\n$leftQuote = [char]0x2018\n$rightQuote = [char]0x2019\n$find2 = $d.Content.Find\n$find2.Text = \"Use only single quotation marks for example (${leftQuote}Hello World${rightQuote}) in the JavaScript code as it can cause errors\"\n$find2.Replacement.Text = \"Use only single quotation marks for example (${leftQuote}Hello World${rightQuote}) in the JavaScript code as it can cause escape errors\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# --- 3. Append a new, mostly-empty row (convention #14) to the table ----\n$table = $d.Tables.Item(1)\n$table.Rows.Add() | Out-Null\n$rowIndex = $table.Rows.Count\n$table.Cell($rowIndex, 1).Range.Text = \"14\"\n\n\"done\"\n"}
